# Update data from Streamlit app
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4 - PT CP Prima
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 45504
$ws.Range("C4").Value = 45869
$ws.Range("F4").Value = 235000000
$ws.Range("G4").Value = 235000000

# ---------------------------------------------------------------------------
# Row 5 - PT Borwita
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 45152
$ws.Range("C5").Value = 45883
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 215000000

# ---------------------------------------------------------------------------
# Row 6 - Shopee Parkiran
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = 45528
$ws.Range("C6").Value = 45893
$ws.Range("G6").Value = 250000000

# ---------------------------------------------------------------------------
# Row 7 - Shopee Express
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = 45572
$ws.Range("C7").Value = 45937
$ws.Range("E7").Value = 1
$ws.Range("G7").Value = 375000000

# ---------------------------------------------------------------------------
# Row 8 - PT Prima Tunggal Mandiri (Shell)
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = 44843
$ws.Range("C8").Value = 45939
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 160000000
$ws.Range("G8").Value = 160000000
$ws.Range("H8").Value = "Split Per Year"

# ---------------------------------------------------------------------------
# Row 9 - PT Prima Tunggal Mandiri (Shell) -- unchanged
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 10 - Oppo Service Center
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 44726
$ws.Range("C10").Value = 45822
$ws.Range("D10").Value = "Reminder: Lease Ending Soon"
$ws.Range("F10").Value = 75000000
$ws.Range("G10").Value = 75000000

# ---------------------------------------------------------------------------
# Row 11 - PT Sukses Muti Servis (Infinix)
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = 45260
$ws.Range("C11").Value = 45991
$ws.Range("E11").Value = 2
$ws.Range("G11").Value = 75000000

# ---------------------------------------------------------------------------
# Row 12 - PT Tumbakmas Niaga Sakti (Sasa)
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = 45291
$ws.Range("C12").Value = 46022
$ws.Range("G12").Value = 525000000

# ---------------------------------------------------------------------------
# Row 13 - PT Guna Elektrik Terang (Philips)
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = 45657
$ws.Range("C13").Value = 46022
$ws.Range("E13").Value = 1
$ws.Range("G13").Value = 180000000
$ws.Range("H13").Value = "Full Lease Upfront"

# ---------------------------------------------------------------------------
# Row 14 (new) - PT Mandiri Akur Pratama (Jogja)
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "PT Mandiri Akur Pratama (Jogja)"
$ws.Range("B14").Value = 45339
$ws.Range("C14").Value = 46070
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 190000000
$ws.Range("G14").Value = 190000000
$ws.Range("H14").Value = "Full Lease Upfront"

# ---------------------------------------------------------------------------
# Row 15 (new) - PT Tiki JNE
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "PT Tiki JNE "
$ws.Range("B15").Value = 45108
$ws.Range("C15").Value = 46569
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 195000000
$ws.Range("G15").Value = 195000000
$ws.Range("H15").Value = "Full Lease Upfront"

# ---------------------------------------------------------------------------
# Row 16 (new) - PT Mandiri Utama Finance
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "PT Mandiri Utama Finance"
$ws.Range("B16").Value = 44166
$ws.Range("C16").Value = 45992
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 235000000
$ws.Range("G16").Value = 235000000
$ws.Range("H16").Value = "Custom Split"
$ws.Range("I16").Value = "50/50/0/0/0"

# ---------------------------------------------------------------------------
# Row 17 (new) - PT Kawan Lama Solusi (Krisbow)
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "PT Kawan Lama Solusi (Krisbow)"
$ws.Range("B17").Value = 45624
$ws.Range("C17").Value = 47450
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 213750000
$ws.Range("G17").Value = 213750000
$ws.Range("H17").Value = "Custom Split"
$ws.Range("I17").Value = "60/0/40/0/0"

# Apply the same date number format (style index 2) used by the other
# Start/Lease End Date cells to the newly added rows, matching the B/C
# columns formatting.
$dateFormat = $ws.Range("B13").NumberFormat
$ws.Range("B14:C17").NumberFormat = $dateFormat
